$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.157.12"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "2.569.96"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'587.05"
$ws.Range("E5").Value = "  +3.28%  "
$ws.Range("D6").Value = "'148.18"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.598"
$ws.Range("E8").Value = "  +1.72%  "
$ws.Range("E9").Value = "  +3.12%  "
$ws.Range("D10").Value = "'5.66"
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").Value = "'27.66"
$ws.Range("D14").Value = "3.029.41"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").Value = "63.102.71"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "'0.0000148"
$ws.Range("E16").Value = "  +2.86%  "
$ws.Range("D17").Value = "2.578.98"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").Value = "'11.38"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'4.45"
$ws.Range("E19").Value = "  +2.90%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'343.26"
$ws.Range("E20").Value = "  +2.20%  "
$ws.Range("E21").Value = "  +1.24%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "'5.54"
$ws.Range("E23").Value = "  -3.54%  "
$ws.Range("D24").Value = "'66.64"
$ws.Range("E24").Value = "  +2.07%  "
$ws.Range("D25").Value = "2.664.93"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("E28").Value = "  +11.42%  "
$ws.Range("D29").Value = "'1.50"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("D31").Value = "'8.44"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").Value = "'1.98"
$ws.Range("E32").Value = "  +6.13%  "
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("D34").Value = "'462.99"
$ws.Range("E34").Value = "  +12.20%  "
$ws.Range("D35").Value = "'176.66"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  +2.89%  "
$ws.Range("D37").Value = "'0.407"
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("D38").Value = "'19.25"
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("E39").Value = "  +4.88%  "
$ws.Range("D41").Value = "'1.76"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "'151.30"
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("D45").Value = "'21.07"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("E46").Value = "  +5.10%  "
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("D48").Value = "'0.0975"
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("D51").Value = "'11.41"
$ws.Range("E51").Value = "  +0.57%  "
